$wb = $excel.ActiveWorkbook

# Update both "展览" and "全部类型" sheets with the same new values in column F
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F5").Value = 3075
    $ws.Range("F7").Value = 2364
    $ws.Range("F11").Value = 1156
    $ws.Range("F15").Value = 1049
    $ws.Range("F16").Value = 285
    $ws.Range("F22").Value = 77
    $ws.Range("F23").Value = 3
    $ws.Range("F24").Value = 17
}
